{"js": "// Build site at 2022-01-09 00:29:46 UTC\n// Update the \"Ativa\u00e7\u00e3o\" date and drop the 4th bullet item (\"heat exchangers\")\n// from the short and long course-program paragraphs (PT + EN).\n\nconst body = context.document.body;\n\nasync function replaceOnce(beforeText, afterText) {\n  const results = body.search(beforeText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${beforeText}`);\n  }\n\n  results.items[0].insertText(afterText, \"Replace\");\n  await context.sync();\n}\n\n// 1) Activation date: 2020 -> 2022\nawait replaceOnce(\n  \"Ativa\u00e7\u00e3o: 01/01/2020\",\n  \"Ativa\u00e7\u00e3o: 01/01/2022\"\n);\n\n// 2) \"Programa resumido\" (PT) \u2014 remove trailing \"4) Trocadores de calor tubulares.\"\nawait replaceOnce(\n  \"Perfis de temperaturas em barras de se\u00e7\u00e3o circular; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido; 4) Trocadores de calor tubulares.\",\n  \"Perfis de temperaturas em barras de se\u00e7\u00e3o circular; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido;\"\n);\n\n// 3) \"Programa resumido\" (EN) \u2014 remove trailing \"4) Tubular heat exchangers.\"\nawait replaceOnce(\n  \"1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems; 4) Tubular heat exchangers.\",\n  \"1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems;\"\n);\n\n// 4) \"Programa\" (PT) \u2014 remove trailing item 4 about heat exchangers\nawait replaceOnce(\n  \"1) Perfis de temperaturas em barras de se\u00e7\u00e3o circular: processos envolvendo condu\u00e7\u00e3o e convec\u00e7\u00e3o em barras de v\u00e1rios materiais e diferentes dimens\u00f5es. Aplica\u00e7\u00e3o do princ\u00edpio das aletas; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o: medidas da varia\u00e7\u00e3o de temperatura em corpos de v\u00e1rias geometrias e materiais diferentes e compara\u00e7\u00e3o com a an\u00e1lise concentrada para regime transiente; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido: avalia\u00e7\u00e3o da transfer\u00eancia de massa entre ar e l\u00edquidos empregando tubos horizontais (c\u00e9lula de Stefan) em regime pseudo-estacion\u00e1rio; 4) Determina\u00e7\u00e3o de coeficientes globais de troca de calor, balan\u00e7os materiais e energ\u00e9ticos em trocadores tubulares do tipo casco e tubos.\",\n  \"1) Perfis de temperaturas em barras de se\u00e7\u00e3o circular: processos envolvendo condu\u00e7\u00e3o e convec\u00e7\u00e3o em barras de v\u00e1rios materiais e diferentes dimens\u00f5es. Aplica\u00e7\u00e3o do princ\u00edpio das aletas; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o: medidas da varia\u00e7\u00e3o de temperatura em corpos de v\u00e1rias geometrias e materiais diferentes e compara\u00e7\u00e3o com a an\u00e1lise concentrada para regime transiente; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido: avalia\u00e7\u00e3o da transfer\u00eancia de massa entre ar e l\u00edquidos empregando tubos horizontais (c\u00e9lula de Stefan) em regime pseudo-estacion\u00e1rio;\"\n);\n\n// 5) \"Programa\" (EN) \u2014 remove trailing item 4, and fix the double space before\n//    \"coefficient\" (\"Diffusion  coefficient\" -> \"Diffusion coefficient\")\nawait replaceOnce(\n  \"1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion  coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state; 4) Determination of overall heat transfer coefficients, material and energetic balances in shell-and-tube heat exchangers.\",\n  \"1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state;\"\n);\n", "ps1": "# Build site at 2022-01-09 00:29:46 UTC\n# Update the \"Ativa\u00e7\u00e3o\" date and drop the 4th bullet item (\"heat exchangers\")\n# from the short and long course-program paragraphs (PT + EN).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $found = $range.Find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n# 1) Activation date: 2020 -> 2022\nReplace-Text \"Ativa\u00e7\u00e3o: 01/01/2020\" \"Ativa\u00e7\u00e3o: 01/01/2022\"\n\n# 2) \"Programa resumido\" (PT) \u2014 remove trailing \"4) Trocadores de calor tubulares.\"\nReplace-Text \"Perfis de temperaturas em barras de se\u00e7\u00e3o circular; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido; 4) Trocadores de calor tubulares.\" \"Perfis de temperaturas em barras de se\u00e7\u00e3o circular; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido;\"\n\n# 3) \"Programa resumido\" (EN) \u2014 remove trailing \"4) Tubular heat exchangers.\"\nReplace-Text \"1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems; 4) Tubular heat exchangers.\" \"1) Temperature distribution in a bar with circular section; 2) Convective heat transfer; 3) Diffusion coefficient in gas-liquid systems;\"\n\n# 4) \"Programa\" (PT) \u2014 remove trailing item 4 about heat exchangers\nReplace-Text \"1) Perfis de temperaturas em barras de se\u00e7\u00e3o circular: processos envolvendo condu\u00e7\u00e3o e convec\u00e7\u00e3o em barras de v\u00e1rios materiais e diferentes dimens\u00f5es. Aplica\u00e7\u00e3o do princ\u00edpio das aletas; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o: medidas da varia\u00e7\u00e3o de temperatura em corpos de v\u00e1rias geometrias e materiais diferentes e compara\u00e7\u00e3o com a an\u00e1lise concentrada para regime transiente; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido: avalia\u00e7\u00e3o da transfer\u00eancia de massa entre ar e l\u00edquidos empregando tubos horizontais (c\u00e9lula de Stefan) em regime pseudo-estacion\u00e1rio; 4) Determina\u00e7\u00e3o de coeficientes globais de troca de calor, balan\u00e7os materiais e energ\u00e9ticos em trocadores tubulares do tipo casco e tubos.\" \"1) Perfis de temperaturas em barras de se\u00e7\u00e3o circular: processos envolvendo condu\u00e7\u00e3o e convec\u00e7\u00e3o em barras de v\u00e1rios materiais e diferentes dimens\u00f5es. Aplica\u00e7\u00e3o do princ\u00edpio das aletas; 2) Transfer\u00eancia de calor por convec\u00e7\u00e3o: medidas da varia\u00e7\u00e3o de temperatura em corpos de v\u00e1rias geometrias e materiais diferentes e compara\u00e7\u00e3o com a an\u00e1lise concentrada para regime transiente; 3) Determina\u00e7\u00e3o do coeficiente de difus\u00e3o em sistemas g\u00e1s-l\u00edquido: avalia\u00e7\u00e3o da transfer\u00eancia de massa entre ar e l\u00edquidos empregando tubos horizontais (c\u00e9lula de Stefan) em regime pseudo-estacion\u00e1rio;\"\n\n# 5) \"Programa\" (EN) \u2014 remove trailing item 4, and fix the double space before\n#    \"coefficient\" (\"Diffusion  coefficient\" -> \"Diffusion coefficient\")\nReplace-Text \"1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion  coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state; 4) Determination of overall heat transfer coefficients, material and energetic balances in shell-and-tube heat exchangers.\" \"1) Temperature distribution in a bar with circular section: heat transfer by conduction and convection in bars of different diameters and materials;; 2) Convective heat transfer: measures temperature variation in bodies of different geometries and materials. Comparison between the experimental data with mathematical models based on the analysis concentrated to transient parameter settings; 3) Diffusion coefficient in gas-liquid systems: analysis of mass transfer between air and liquids using horizontal pipes (Stefan cell) in pseudo-steady state;\"\n"}
